$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix provider/name fields: commas -> periods (plain text, unambiguous) ---
$ws.Range('E62').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E88').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('F88').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('E94').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F94').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E95').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E125').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('F125').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('E135').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('F142').Value = 'MERCANZINI. GASTON ARIEL'
$ws.Range('E156').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E166').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E185').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E204').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E221').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'

# --- Fix "Importe" amount fields: Spanish-formatted numeric text (e.g. 1.234,56) -> plain
# decimal-point text (1234.56). These cells must stay TEXT (shared-string) cells, not turn
# into real numbers, so the value is entered with a leading apostrophe (forces text entry,
# same as typing it by hand in Excel) and the cell style is then put back to "Normal" so no
# extra formatting is left behind on the cell itself.
$c = $ws.Range('H2')
$c.Value = "'2979.00"
$c.Style = "Normal"
$c = $ws.Range('H3')
$c.Value = "'176000.00"
$c.Style = "Normal"
$c = $ws.Range('H4')
$c.Value = "'163550.00"
$c.Style = "Normal"
$c = $ws.Range('H5')
$c.Value = "'20000.00"
$c.Style = "Normal"
$c = $ws.Range('H6')
$c.Value = "'50000.00"
$c.Style = "Normal"
$c = $ws.Range('H7')
$c.Value = "'817852.00"
$c.Style = "Normal"
$c = $ws.Range('H8')
$c.Value = "'470.00"
$c.Style = "Normal"
$c = $ws.Range('H9')
$c.Value = "'3665.00"
$c.Style = "Normal"
$c = $ws.Range('H10')
$c.Value = "'12900.00"
$c.Style = "Normal"
$c = $ws.Range('H11')
$c.Value = "'414837.57"
$c.Style = "Normal"
$c = $ws.Range('H12')
$c.Value = "'176.00"
$c.Style = "Normal"
$c = $ws.Range('H13')
$c.Value = "'726000.00"
$c.Style = "Normal"
$c = $ws.Range('H14')
$c.Value = "'374700.00"
$c.Style = "Normal"
$c = $ws.Range('H15')
$c.Value = "'473.60"
$c.Style = "Normal"
$c = $ws.Range('H16')
$c.Value = "'1050.00"
$c.Style = "Normal"
$c = $ws.Range('H17')
$c.Value = "'260.00"
$c.Style = "Normal"
$c = $ws.Range('H18')
$c.Value = "'169802.93"
$c.Style = "Normal"
$c = $ws.Range('H19')
$c.Value = "'570040.04"
$c.Style = "Normal"
$c = $ws.Range('H20')
$c.Value = "'84227.95"
$c.Style = "Normal"
$c = $ws.Range('H21')
$c.Value = "'9400.00"
$c.Style = "Normal"
$c = $ws.Range('H22')
$c.Value = "'1050.00"
$c.Style = "Normal"
$c = $ws.Range('H23')
$c.Value = "'3192.00"
$c.Style = "Normal"
$c = $ws.Range('H24')
$c.Value = "'11834.27"
$c.Style = "Normal"
$c = $ws.Range('H25')
$c.Value = "'18546.32"
$c.Style = "Normal"
$c = $ws.Range('H26')
$c.Value = "'14715.00"
$c.Style = "Normal"
$c = $ws.Range('H27')
$c.Value = "'646.40"
$c.Style = "Normal"
$c = $ws.Range('H28')
$c.Value = "'3200.00"
$c.Style = "Normal"
$c = $ws.Range('H29')
$c.Value = "'140.48"
$c.Style = "Normal"
$c = $ws.Range('H30')
$c.Value = "'4810.00"
$c.Style = "Normal"
$c = $ws.Range('H31')
$c.Value = "'1200.00"
$c.Style = "Normal"
$c = $ws.Range('H32')
$c.Value = "'121.92"
$c.Style = "Normal"
$c = $ws.Range('H33')
$c.Value = "'1958.74"
$c.Style = "Normal"
$c = $ws.Range('H34')
$c.Value = "'74.00"
$c.Style = "Normal"
$c = $ws.Range('H35')
$c.Value = "'51403.81"
$c.Style = "Normal"
$c = $ws.Range('H36')
$c.Value = "'240.00"
$c.Style = "Normal"
$c = $ws.Range('H37')
$c.Value = "'165.00"
$c.Style = "Normal"
$c = $ws.Range('H38')
$c.Value = "'500.00"
$c.Style = "Normal"
$c = $ws.Range('H39')
$c.Value = "'71.83"
$c.Style = "Normal"
$c = $ws.Range('H40')
$c.Value = "'116853.00"
$c.Style = "Normal"
$c = $ws.Range('H41')
$c.Value = "'12828.40"
$c.Style = "Normal"
$c = $ws.Range('H42')
$c.Value = "'424.88"
$c.Style = "Normal"
$c = $ws.Range('H43')
$c.Value = "'1299.00"
$c.Style = "Normal"
$c = $ws.Range('H44')
$c.Value = "'1478.00"
$c.Style = "Normal"
$c = $ws.Range('H45')
$c.Value = "'934.36"
$c.Style = "Normal"
$c = $ws.Range('H46')
$c.Value = "'14166.00"
$c.Style = "Normal"
$c = $ws.Range('H47')
$c.Value = "'7565.00"
$c.Style = "Normal"
$c = $ws.Range('H48')
$c.Value = "'58851.21"
$c.Style = "Normal"
$c = $ws.Range('H49')
$c.Value = "'2499.66"
$c.Style = "Normal"
$c = $ws.Range('H50')
$c.Value = "'751.18"
$c.Style = "Normal"
$c = $ws.Range('H51')
$c.Value = "'3736.03"
$c.Style = "Normal"
$c = $ws.Range('H52')
$c.Value = "'10087.00"
$c.Style = "Normal"
$c = $ws.Range('H53')
$c.Value = "'28.75"
$c.Style = "Normal"
$c = $ws.Range('H54')
$c.Value = "'1750.53"
$c.Style = "Normal"
$c = $ws.Range('H55')
$c.Value = "'396.00"
$c.Style = "Normal"
$c = $ws.Range('H56')
$c.Value = "'810.00"
$c.Style = "Normal"
$c = $ws.Range('H57')
$c.Value = "'877.58"
$c.Style = "Normal"
$c = $ws.Range('H58')
$c.Value = "'1219.00"
$c.Style = "Normal"
$c = $ws.Range('H59')
$c.Value = "'2480.00"
$c.Style = "Normal"
$c = $ws.Range('H60')
$c.Value = "'42354.09"
$c.Style = "Normal"
$c = $ws.Range('H61')
$c.Value = "'330.00"
$c.Style = "Normal"
$c = $ws.Range('H62')
$c.Value = "'530.00"
$c.Style = "Normal"
$c = $ws.Range('H63')
$c.Value = "'32254.50"
$c.Style = "Normal"
$c = $ws.Range('H64')
$c.Value = "'359.04"
$c.Style = "Normal"
$c = $ws.Range('H65')
$c.Value = "'1070.99"
$c.Style = "Normal"
$c = $ws.Range('H66')
$c.Value = "'4646.24"
$c.Style = "Normal"
$c = $ws.Range('H67')
$c.Value = "'962.54"
$c.Style = "Normal"
$c = $ws.Range('H68')
$c.Value = "'1555.00"
$c.Style = "Normal"
$c = $ws.Range('H69')
$c.Value = "'805.00"
$c.Style = "Normal"
$c = $ws.Range('H70')
$c.Value = "'240.00"
$c.Style = "Normal"
$c = $ws.Range('H71')
$c.Value = "'22631.68"
$c.Style = "Normal"
$c = $ws.Range('H72')
$c.Value = "'765.77"
$c.Style = "Normal"
$c = $ws.Range('H73')
$c.Value = "'329.00"
$c.Style = "Normal"
$c = $ws.Range('H74')
$c.Value = "'793.71"
$c.Style = "Normal"
$c = $ws.Range('H75')
$c.Value = "'270.00"
$c.Style = "Normal"
$c = $ws.Range('H76')
$c.Value = "'1026.00"
$c.Style = "Normal"
$c = $ws.Range('H77')
$c.Value = "'2040.65"
$c.Style = "Normal"
$c = $ws.Range('H78')
$c.Value = "'150.00"
$c.Style = "Normal"
$c = $ws.Range('H79')
$c.Value = "'17834.00"
$c.Style = "Normal"
$c = $ws.Range('H80')
$c.Value = "'2875.00"
$c.Style = "Normal"
$c = $ws.Range('H81')
$c.Value = "'5976.12"
$c.Style = "Normal"
$c = $ws.Range('H82')
$c.Value = "'1570.00"
$c.Style = "Normal"
$c = $ws.Range('H83')
$c.Value = "'520.00"
$c.Style = "Normal"
$c = $ws.Range('H84')
$c.Value = "'31800.00"
$c.Style = "Normal"
$c = $ws.Range('H85')
$c.Value = "'2620.00"
$c.Style = "Normal"
$c = $ws.Range('H86')
$c.Value = "'646.00"
$c.Style = "Normal"
$c = $ws.Range('H87')
$c.Value = "'14856.00"
$c.Style = "Normal"
$c = $ws.Range('H88')
$c.Value = "'60.00"
$c.Style = "Normal"
$c = $ws.Range('H89')
$c.Value = "'2292.00"
$c.Style = "Normal"
$c = $ws.Range('H90')
$c.Value = "'5745.00"
$c.Style = "Normal"
$c = $ws.Range('H91')
$c.Value = "'81.80"
$c.Style = "Normal"
$c = $ws.Range('H92')
$c.Value = "'558.00"
$c.Style = "Normal"
$c = $ws.Range('H93')
$c.Value = "'504.00"
$c.Style = "Normal"
$c = $ws.Range('H94')
$c.Value = "'763.63"
$c.Style = "Normal"
$c = $ws.Range('H95')
$c.Value = "'1890.00"
$c.Style = "Normal"
$c = $ws.Range('H96')
$c.Value = "'591.30"
$c.Style = "Normal"
$c = $ws.Range('H97')
$c.Value = "'4355.88"
$c.Style = "Normal"
$c = $ws.Range('H98')
$c.Value = "'1520.00"
$c.Style = "Normal"
$c = $ws.Range('H99')
$c.Value = "'44000.00"
$c.Style = "Normal"
$c = $ws.Range('H100')
$c.Value = "'4500.00"
$c.Style = "Normal"
$c = $ws.Range('H101')
$c.Value = "'5500.00"
$c.Style = "Normal"
$c = $ws.Range('H102')
$c.Value = "'7000.00"
$c.Style = "Normal"
$c = $ws.Range('H103')
$c.Value = "'45000.00"
$c.Style = "Normal"
$c = $ws.Range('H104')
$c.Value = "'6000.00"
$c.Style = "Normal"
$c = $ws.Range('H105')
$c.Value = "'3500.00"
$c.Style = "Normal"
$c = $ws.Range('H106')
$c.Value = "'5000.00"
$c.Style = "Normal"
$c = $ws.Range('H107')
$c.Value = "'23000.00"
$c.Style = "Normal"
$c = $ws.Range('H108')
$c.Value = "'60000.00"
$c.Style = "Normal"
$c = $ws.Range('H109')
$c.Value = "'28500.00"
$c.Style = "Normal"
$c = $ws.Range('H110')
$c.Value = "'40500.00"
$c.Style = "Normal"
$c = $ws.Range('H111')
$c.Value = "'162.58"
$c.Style = "Normal"
$c = $ws.Range('H112')
$c.Value = "'500.00"
$c.Style = "Normal"
$c = $ws.Range('H113')
$c.Value = "'270000.00"
$c.Style = "Normal"
$c = $ws.Range('H114')
$c.Value = "'9469.51"
$c.Style = "Normal"
$c = $ws.Range('H115')
$c.Value = "'880.32"
$c.Style = "Normal"
$c = $ws.Range('H116')
$c.Value = "'63823.54"
$c.Style = "Normal"
$c = $ws.Range('H117')
$c.Value = "'120.69"
$c.Style = "Normal"
$c = $ws.Range('H118')
$c.Value = "'56.10"
$c.Style = "Normal"
$c = $ws.Range('H119')
$c.Value = "'13650.00"
$c.Style = "Normal"
$c = $ws.Range('H120')
$c.Value = "'14818.80"
$c.Style = "Normal"
$c = $ws.Range('H121')
$c.Value = "'300.00"
$c.Style = "Normal"
$c = $ws.Range('H122')
$c.Value = "'2896.30"
$c.Style = "Normal"
$c = $ws.Range('H123')
$c.Value = "'28.42"
$c.Style = "Normal"
$c = $ws.Range('H124')
$c.Value = "'180.00"
$c.Style = "Normal"
$c = $ws.Range('H125')
$c.Value = "'250.00"
$c.Style = "Normal"
$c = $ws.Range('H126')
$c.Value = "'1840.00"
$c.Style = "Normal"
$c = $ws.Range('H127')
$c.Value = "'6978.48"
$c.Style = "Normal"
$c = $ws.Range('H128')
$c.Value = "'307.30"
$c.Style = "Normal"
$c = $ws.Range('H129')
$c.Value = "'11093.32"
$c.Style = "Normal"
$c = $ws.Range('H130')
$c.Value = "'415.00"
$c.Style = "Normal"
$c = $ws.Range('H131')
$c.Value = "'38.94"
$c.Style = "Normal"
$c = $ws.Range('H132')
$c.Value = "'1450.00"
$c.Style = "Normal"
$c = $ws.Range('H133')
$c.Value = "'37.29"
$c.Style = "Normal"
$c = $ws.Range('H134')
$c.Value = "'22.76"
$c.Style = "Normal"
$c = $ws.Range('H135')
$c.Value = "'260.00"
$c.Style = "Normal"
$c = $ws.Range('H136')
$c.Value = "'14200.00"
$c.Style = "Normal"
$c = $ws.Range('H137')
$c.Value = "'1980.00"
$c.Style = "Normal"
$c = $ws.Range('H138')
$c.Value = "'4476.00"
$c.Style = "Normal"
$c = $ws.Range('H139')
$c.Value = "'3120.00"
$c.Style = "Normal"
$c = $ws.Range('H140')
$c.Value = "'5900.00"
$c.Style = "Normal"
$c = $ws.Range('H141')
$c.Value = "'1485.00"
$c.Style = "Normal"
$c = $ws.Range('H142')
$c.Value = "'6000.00"
$c.Style = "Normal"
$c = $ws.Range('H143')
$c.Value = "'8000.00"
$c.Style = "Normal"
$c = $ws.Range('H144')
$c.Value = "'65.85"
$c.Style = "Normal"
$c = $ws.Range('H145')
$c.Value = "'2534.87"
$c.Style = "Normal"
$c = $ws.Range('H146')
$c.Value = "'1033.71"
$c.Style = "Normal"
$c = $ws.Range('H147')
$c.Value = "'3107.95"
$c.Style = "Normal"
$c = $ws.Range('H148')
$c.Value = "'128799.90"
$c.Style = "Normal"
$c = $ws.Range('H149')
$c.Value = "'131598.40"
$c.Style = "Normal"
$c = $ws.Range('H150')
$c.Value = "'2880.00"
$c.Style = "Normal"
$c = $ws.Range('H151')
$c.Value = "'5000.00"
$c.Style = "Normal"
$c = $ws.Range('H152')
$c.Value = "'6630.00"
$c.Style = "Normal"
$c = $ws.Range('H153')
$c.Value = "'950.00"
$c.Style = "Normal"
$c = $ws.Range('H154')
$c.Value = "'290.00"
$c.Style = "Normal"
$c = $ws.Range('H155')
$c.Value = "'4480.00"
$c.Style = "Normal"
$c = $ws.Range('H156')
$c.Value = "'350.00"
$c.Style = "Normal"
$c = $ws.Range('H157')
$c.Value = "'300.00"
$c.Style = "Normal"
$c = $ws.Range('H158')
$c.Value = "'2620.00"
$c.Style = "Normal"
$c = $ws.Range('H159')
$c.Value = "'200.00"
$c.Style = "Normal"
$c = $ws.Range('H160')
$c.Value = "'1074.48"
$c.Style = "Normal"
$c = $ws.Range('H161')
$c.Value = "'25160.00"
$c.Style = "Normal"
$c = $ws.Range('H162')
$c.Value = "'1300.00"
$c.Style = "Normal"
$c = $ws.Range('H163')
$c.Value = "'480.00"
$c.Style = "Normal"
$c = $ws.Range('H164')
$c.Value = "'10836.28"
$c.Style = "Normal"
$c = $ws.Range('H165')
$c.Value = "'150.00"
$c.Style = "Normal"
$c = $ws.Range('H166')
$c.Value = "'1050.00"
$c.Style = "Normal"
$c = $ws.Range('H167')
$c.Value = "'99.50"
$c.Style = "Normal"
$c = $ws.Range('H168')
$c.Value = "'2713.56"
$c.Style = "Normal"
$c = $ws.Range('H169')
$c.Value = "'1000.00"
$c.Style = "Normal"
$c = $ws.Range('H170')
$c.Value = "'2814.64"
$c.Style = "Normal"
$c = $ws.Range('H171')
$c.Value = "'595.00"
$c.Style = "Normal"
$c = $ws.Range('H172')
$c.Value = "'3305.00"
$c.Style = "Normal"
$c = $ws.Range('H173')
$c.Value = "'1710.00"
$c.Style = "Normal"
$c = $ws.Range('H174')
$c.Value = "'826.79"
$c.Style = "Normal"
$c = $ws.Range('H175')
$c.Value = "'909.00"
$c.Style = "Normal"
$c = $ws.Range('H176')
$c.Value = "'3621.92"
$c.Style = "Normal"
$c = $ws.Range('H177')
$c.Value = "'2496.08"
$c.Style = "Normal"
$c = $ws.Range('H178')
$c.Value = "'377.40"
$c.Style = "Normal"
$c = $ws.Range('H179')
$c.Value = "'706.92"
$c.Style = "Normal"
$c = $ws.Range('H180')
$c.Value = "'432.35"
$c.Style = "Normal"
$c = $ws.Range('H181')
$c.Value = "'18560.40"
$c.Style = "Normal"
$c = $ws.Range('H182')
$c.Value = "'220.00"
$c.Style = "Normal"
$c = $ws.Range('H183')
$c.Value = "'8956.84"
$c.Style = "Normal"
$c = $ws.Range('H184')
$c.Value = "'4000.00"
$c.Style = "Normal"
$c = $ws.Range('H185')
$c.Value = "'6311.00"
$c.Style = "Normal"
$c = $ws.Range('H186')
$c.Value = "'3482.76"
$c.Style = "Normal"
$c = $ws.Range('H187')
$c.Value = "'24.14"
$c.Style = "Normal"
$c = $ws.Range('H188')
$c.Value = "'23936.00"
$c.Style = "Normal"
$c = $ws.Range('H189')
$c.Value = "'496.00"
$c.Style = "Normal"
$c = $ws.Range('H190')
$c.Value = "'110.00"
$c.Style = "Normal"
$c = $ws.Range('H191')
$c.Value = "'1770.00"
$c.Style = "Normal"
$c = $ws.Range('H192')
$c.Value = "'123.00"
$c.Style = "Normal"
$c = $ws.Range('H193')
$c.Value = "'1290.25"
$c.Style = "Normal"
$c = $ws.Range('H194')
$c.Value = "'1527.60"
$c.Style = "Normal"
$c = $ws.Range('H195')
$c.Value = "'4699.23"
$c.Style = "Normal"
$c = $ws.Range('H196')
$c.Value = "'300.00"
$c.Style = "Normal"
$c = $ws.Range('H197')
$c.Value = "'2050.00"
$c.Style = "Normal"
$c = $ws.Range('H198')
$c.Value = "'2887.50"
$c.Style = "Normal"
$c = $ws.Range('H199')
$c.Value = "'21000.00"
$c.Style = "Normal"
$c = $ws.Range('H200')
$c.Value = "'1962056.83"
$c.Style = "Normal"
$c = $ws.Range('H201')
$c.Value = "'18500.00"
$c.Style = "Normal"
$c = $ws.Range('H202')
$c.Value = "'4520.84"
$c.Style = "Normal"
$c = $ws.Range('H203')
$c.Value = "'580.00"
$c.Style = "Normal"
$c = $ws.Range('H204')
$c.Value = "'280.00"
$c.Style = "Normal"
$c = $ws.Range('H205')
$c.Value = "'3514.79"
$c.Style = "Normal"
$c = $ws.Range('H206')
$c.Value = "'567.73"
$c.Style = "Normal"
$c = $ws.Range('H207')
$c.Value = "'2900.00"
$c.Style = "Normal"
$c = $ws.Range('H208')
$c.Value = "'433000.00"
$c.Style = "Normal"
$c = $ws.Range('H209')
$c.Value = "'100360.00"
$c.Style = "Normal"
$c = $ws.Range('H210')
$c.Value = "'65500.00"
$c.Style = "Normal"
$c = $ws.Range('H211')
$c.Value = "'319321.64"
$c.Style = "Normal"
$c = $ws.Range('H212')
$c.Value = "'60000.00"
$c.Style = "Normal"
$c = $ws.Range('H213')
$c.Value = "'40000.00"
$c.Style = "Normal"
$c = $ws.Range('H214')
$c.Value = "'115500.00"
$c.Style = "Normal"
$c = $ws.Range('H215')
$c.Value = "'411312.00"
$c.Style = "Normal"
$c = $ws.Range('H216')
$c.Value = "'223000.00"
$c.Style = "Normal"
$c = $ws.Range('H217')
$c.Value = "'446000.00"
$c.Style = "Normal"
$c = $ws.Range('H218')
$c.Value = "'432000.00"
$c.Style = "Normal"
$c = $ws.Range('H219')
$c.Value = "'65500.00"
$c.Style = "Normal"
$c = $ws.Range('H220')
$c.Value = "'104000.00"
$c.Style = "Normal"
$c = $ws.Range('H221')
$c.Value = "'480.00"
$c.Style = "Normal"
$c = $ws.Range('H222')
$c.Value = "'303471.00"
$c.Style = "Normal"
$c = $ws.Range('H223')
$c.Value = "'20720.00"
$c.Style = "Normal"
$c = $ws.Range('H224')
$c.Value = "'9000.00"
$c.Style = "Normal"
$c = $ws.Range('H225')
$c.Value = "'18000.00"
$c.Style = "Normal"
$c = $ws.Range('H226')
$c.Value = "'5000.00"
$c.Style = "Normal"
$c = $ws.Range('H227')
$c.Value = "'1110.00"
$c.Style = "Normal"
